# Applies the "operazioni" worksheet edit:
#  - Updates several "Frequenza" (column C) values
#  - Removes the image-detail operation (old O8 row) and shifts the
#    remaining rows up so O9-O13 become O8-O11 (text-wise) while the
#    O-codes in column A stay the same for every surviving row
#  - Deletes the two trailing rows (old O12/O13) entirely
#  - Adjusts row 9's height and column C's width to match the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Frequenza) textual updates for rows that keep their code ---
$ws.Range("C2").Value = "1 all’anno"
$ws.Range("C4").Value = "1 ogni 2 settimane"
$ws.Range("C5").Value = "10 ogni 2 settimane"
$ws.Range("C6").Value = "1 ogni 2 anni"

# --- Shift the descriptions: O8's old (image-related) description is
#     dropped, and every following row's text moves up by one slot ---
$ws.Range("B9").Value  = "Visualizzare gli organismi avvistati in una  spedizione"
$ws.Range("C9").Value  = "1 al giorno"

$ws.Range("B10").Value = "Visualizzare le formazioni geologiche di un determinato grado di pericolo e dove sono situate"
$ws.Range("C10").Value = "1 al mese"

$ws.Range("B11").Value = "Visualizzare i luoghi dove sono affondati determinati relitti (sapendo il nome)"
$ws.Range("C11").Value = "1 al mese"

$ws.Range("B12").Value = "Visualizzare le analisi fatte su un materiale e da quale laboratorio sono state eseguite"
$ws.Range("C12").Value = "1 al giorno"

# --- Remove the two now-obsolete trailing rows (old O12 / O13, which
#     described the removed image-backed operations) ---
$ws.Rows("13:14").Delete()

# --- Row 9 shrinks (shorter text, no longer needs the tall wrap) ---
$ws.Rows("9:9").RowHeight = 31.3

# --- Column C widens slightly to fit "10 ogni 2 settimane" ---
$ws.Columns("C:C").ColumnWidth = 20

$ws.Range("C16").Select() | Out-Null
